$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.21875
$ws.Range("C2").Value = 0.5
$ws.Range("J2").Value = 0.0625
$ws.Range("P2").Value = 0.15625
$ws.Range("S2").Value = 0.0625
$ws.Range("J3").Value = 0.1333333333333333
$ws.Range("P3").Value = 0.6
$ws.Range("S3").Value = 0.2666666666666667
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("F6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.2666666666666667
$ws.Range("Q6").Value = 0.2
$ws.Range("S6").Value = 0.4
$ws.Range("D7").Value = 0.08333333333333333
$ws.Range("F7").Value = 0.08333333333333333
$ws.Range("J7").Value = 0.08333333333333333
$ws.Range("Q7").Value = 0.3333333333333333
$ws.Range("R7").Value = 0.08333333333333333
$ws.Range("S7").Value = 0.3333333333333333
$ws.Range("B8").Value = 0.08333333333333333
$ws.Range("D8").Value = 0.05555555555555555
$ws.Range("J8").Value = 0.1944444444444444
$ws.Range("O8").Value = 0.02777777777777778
$ws.Range("Q8").Value = 0.1944444444444444
$ws.Range("R8").Value = 0.1666666666666667
$ws.Range("S8").Value = 0.2777777777777778
$ws.Range("B9").Value = 0.1428571428571428
$ws.Range("F9").Value = 0.2857142857142857
$ws.Range("J9").Value = 0.2142857142857143
$ws.Range("Q9").Value = 0.1428571428571428
$ws.Range("R9").Value = 0.1428571428571428
$ws.Range("S9").Value = 0.07142857142857142
$ws.Range("B10").Value = 0.1637931034482759
$ws.Range("D10").Value = 0.008620689655172414
$ws.Range("F10").Value = 0.04310344827586207
$ws.Range("J10").Value = 0.1206896551724138
$ws.Range("O10").Value = 0.01724137931034483
$ws.Range("Q10").Value = 0.2413793103448276
$ws.Range("R10").Value = 0.103448275862069
$ws.Range("S10").Value = 0.3017241379310345
$ws.Range("G11").Value = 0.2592592592592592
$ws.Range("J11").Value = 0.1481481481481481
$ws.Range("K11").Value = 0.4074074074074074
$ws.Range("L11").Value = 0.1851851851851852
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.1666666666666667
$ws.Range("K12").Value = 0.1666666666666667
$ws.Range("F13").Value = 0.25
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.25
$ws.Range("H15").Value = 0.2
$ws.Range("O15").Value = 0.1
$ws.Range("S15").Value = 0.2
$ws.Range("H16").Value = 0.2857142857142857
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.5
$ws.Range("K16").Value = 0.07142857142857142
$ws.Range("S16").Value = 0.07142857142857142
$ws.Range("H17").Value = 0.2045454545454546
$ws.Range("I17").Value = 0.04545454545454546
$ws.Range("J17").Value = 0.4090909090909091
$ws.Range("K17").Value = 0.09090909090909091
$ws.Range("M17").Value = 0.04545454545454546
$ws.Range("O17").Value = 0.06818181818181818
$ws.Range("S17").Value = 0.1363636363636364
$ws.Range("H18").Value = 0.2857142857142857
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.09523809523809523
$ws.Range("M18").Value = 0.09523809523809523
$ws.Range("S18").Value = 0.09523809523809523
$ws.Range("H19").Value = 0.186046511627907
$ws.Range("I19").Value = 0.1279069767441861
$ws.Range("J19").Value = 0.4534883720930232
$ws.Range("K19").Value = 0.05813953488372093
$ws.Range("M19").Value = 0.01162790697674419
$ws.Range("O19").Value = 0.02325581395348837
$ws.Range("S19").Value = 0.1395348837209302
